# Apply updated cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.334.40"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.715.05"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.15"
$ws.Range("E5").Value = "  -2.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5307"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2639"
$ws.Range("E8").Value = "  -4.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06548"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.06"
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07653"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.568"
$ws.Range("E12").Value = "  -3.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.956.16"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.687.79"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5746"
$ws.Range("E15").Value = "  -5.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8181"
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.34"
$ws.Range("E17").Value = "  -2.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.357.22"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.22"
$ws.Range("E19").Value = "  +2.77%  "
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.687"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("E22").Value = "  -4.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.959"
$ws.Range("E23").Value = "  -4.19%  "
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.10"
$ws.Range("E25").Value = "  -2.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.734"
$ws.Range("E26").Value = "  +8.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1218"
$ws.Range("E27").Value = "  -3.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.272"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.35"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05411"
$ws.Range("E30").Value = "  -3.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.293"
$ws.Range("E31").Value = "  -2.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.493"
$ws.Range("E32").Value = "  -4.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.409"
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.636"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.874"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.432"
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9496"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5866"
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01635"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.878"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8415"
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.038.50"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.06"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.860.86"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("E46").Value = "  +5.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.12"
$ws.Range("E47").Value = "  -3.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4503"
$ws.Range("E48").Value = "  +3.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06638"
$ws.Range("E49").Value = "  +15.05%  "
$ws.Range("B50").Value = "Frax"
$ws.Range("C50").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.004"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.092"
$ws.Range("E51").Value = "  -0.03%  "
